# Applies the "Updated UG, PPP, DG" edit:
#  1. Re-cache the cover-slide date placeholder ("datetimeFigureOut" field)
#     from 4/10/2019 to 4/14/2019 everywhere it is defined (slide master,
#     every slide layout, and the notes master) - mirrors every slide in
#     the authored deck picking up the new save date.
#  2. Resize/reposition the ":Finance / BookParser" actor box on the
#     sequence diagram and rename it to "TrackerParser".

$p = $ppt.ActivePresentation

$oldDate = "4/10/2019"
$newDate = "4/14/2019"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# -- 1a. Slide master date placeholder --
Update-DatePlaceholder $p.SlideMaster.Shapes

# -- 1b. Every slide layout's date placeholder --
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# -- 1c. Notes master date placeholder --
Update-DatePlaceholder $p.NotesMaster.Shapes

# -- 2. Resize / rename the actor box on slide 1 --
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item(8)

$shape.Left = 198.44086624173227
$shape.Top = -3.7924410448818895
$shape.Width = 111.3759880519685
$shape.Height = 34.55984311968504

$para = $shape.TextFrame.TextRange.Paragraphs(2, 1)
$para.Runs(1, 1).Text = "TrackerParser"
